$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("B2").Value = "d"
$ws.Range("C2").Value = "d"
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = 3
$ws.Range("H2").Value = 51.48

# Row 3 updates
$ws.Range("A3").ClearContents()
$ws.Range("B3").Value = "r"
$ws.Range("C3").Value = "r"
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 5
$ws.Range("H3").Value = 91.16

# Row 4: clear all data cells except A4
$ws.Range("B4:H4").ClearContents()

# Row 5: clear A5
$ws.Range("A5").ClearContents()
